# button1_Click (fix enter failse PPhone) Add Form (Delete & Update)
# Rebuild the member list: clear the old rows and re-enter the data via the
# Add/Update/Delete form, then append the 4 time-slot lookup rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Delete old data -------------------------------------------------
$ws.Cells.Clear()

# --- New member data --------------------------------------------------
$ids    = @(1, 2, 3, 4, 5, 6)
$names  = @("John Smith", "Emily Johnson", "Michael Davis", "Sarah Thompson", "David Brown", "Jennifer Wilson")
$phones = @(12345, 12346, 12347, 12348, 12349, 12341)
$genders= @("Male", "Male", "Female", "Female", "Male", "Male")
$ages   = @(20, 18, 10, 20, 17, 15)
$fees   = @(2000, 1000, 3000, 4000, 5000, 500)
$slots  = @("6AM-8AM", "8AM-10AM", "6PM-8PM", "8PM-10PM", "8PM-10PM", "6PM-8PM")

# Column A : Id
for ($i = 0; $i -lt 6; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $ids[$i]
}

# Column B : Name
for ($i = 0; $i -lt 6; $i++) {
    $ws.Cells.Item($i + 1, 2).Value = $names[$i]
}

# Column C : Phone
for ($i = 0; $i -lt 6; $i++) {
    $ws.Cells.Item($i + 1, 3).Value = $phones[$i]
}

# Column D : Gender
for ($i = 0; $i -lt 6; $i++) {
    $ws.Cells.Item($i + 1, 4).Value = $genders[$i]
}

# Column E : Age
for ($i = 0; $i -lt 6; $i++) {
    $ws.Cells.Item($i + 1, 5).Value = $ages[$i]
}

# Column F : Fee
for ($i = 0; $i -lt 6; $i++) {
    $ws.Cells.Item($i + 1, 6).Value = $fees[$i]
}

# Column G : Time slot
for ($i = 0; $i -lt 6; $i++) {
    $ws.Cells.Item($i + 1, 7).Value = $slots[$i]
}

# --- Time-slot lookup list appended below the table (A7:A10) ---------
$lookupSlots = @("6AM-8AM", "8AM-10AM", "6PM-8PM", "8PM-10PM")
for ($i = 0; $i -lt 4; $i++) {
    $ws.Cells.Item(7 + $i, 1).Value = $lookupSlots[$i]
}

# --- Restore the selection left behind by the form ---------------------
$ws.Range("A3:G4").Select()
